$p = $ppt.ActivePresentation
try { $app.DisplayGuides = $true } catch {}
try { $null = $p.SlideMaster.Guides.Add(1, 3.0) } catch {}
Write-Output "done"
